# Voice Input languages workbook — add the "toggle" / "Toggle" / "切換"
# translation row (row 31) to Tabelle1, matching the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: key / English / Traditional-Chinese translation.
$ws.Range("A31").Value = "toggle"
$ws.Range("B31").Value = "Toggle"
$ws.Range("C31").Value = "切換"

# Column C carries the Chinese-font cell style used throughout the sheet
# (Microsoft JhengHei). Copy formatting from the cell directly above so the
# new row matches instead of minting a brand-new style entry.
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null

# Leave the sheet's selection on the newly-added row, as it was after the
# edit in Excel.
$ws.Range("B34").Select() | Out-Null
